$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the F column values (rows 4-8)
$ws.Range("F4").Value = 99
$ws.Range("F5").Value = 98
$ws.Range("F6").Value = 99
$ws.Range("F7").Value = 99
$ws.Range("F8").Value = 98

# Update the active selection to F9
$ws.Range("F9").Select()
